$d = $word.ActiveDocument

function New-PkgXml($innerBody) {
    $pre = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
    $post = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $pre + $innerBody + $post
}

# ---------------------------------------------------------------------------
# Edit 1: insert a new lead-in paragraph ("At first, I plan to analyse all
# applications including free and paid. However, ") right before the
# "Clean App Store Data" paragraph.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Clean App Store Data", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para1 = $rng1.Paragraphs(1)

$newLead = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">At first, I plan to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>analyse</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> all applications including free and paid. However, </w:t></w:r></w:p>'
$origPara1 = '<w:p w14:paraId="5AD78399" w14:textId="77777777" w:rsidR="007130E0" w:rsidRPr="00FD000B" w:rsidRDefault="007130E0" w:rsidP="007130E0"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00FD000B"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Clean App Store Data</w:t></w:r></w:p>'

$para1.Range.InsertXML((New-PkgXml ($newLead + $origPara1)))

# ---------------------------------------------------------------------------
# Edit 2: split the run in the "Delete Column application size..." bullet so
# that "code,  number" gets its own run wrapped in gramStart/gramEnd proof
# markers.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Column application size, latest version code,  number of supporting devices", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2 = $rng2.Paragraphs(1)

$newPara2 = '<w:p w14:paraId="71F7EDE3" w14:textId="77777777" w:rsidR="007130E0" w:rsidRPr="00FC5624" w:rsidRDefault="007130E0" w:rsidP="007130E0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Delete </w:t></w:r><w:r w:rsidRPr="00FC5624"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Column application size, latest version </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>code,  number</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> of supporting devices, number of screenshots showed for display, number of supported languages, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00FC5624"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Vpp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00FC5624"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> device based licensing enabled, we are not interested in if </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>those features</w:t></w:r><w:r w:rsidRPr="00FC5624"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> will affect number of installs and rating</w:t></w:r></w:p>'

$para2.Range.InsertXML((New-PkgXml $newPara2))

# ---------------------------------------------------------------------------
# Edit 3: split the run in "How many category for rating" so "category" gets
# its own run wrapped in gramStart/gramEnd proof markers.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("How many category for rating", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para3 = $rng3.Paragraphs(1)

$newPara3 = '<w:p w14:paraId="7D863C08" w14:textId="77777777" w:rsidR="007130E0" w:rsidRDefault="007130E0" w:rsidP="007130E0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">How many </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>category</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> for rating </w:t></w:r></w:p>'

$para3.Range.InsertXML((New-PkgXml $newPara3))

# ---------------------------------------------------------------------------
# Edit 4: add a <w:lastRenderedPageBreak/> before the text run that starts
# "Delete number of language support column in App ".
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("Delete number of language support column in App", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para4 = $rng4.Paragraphs(1)

$newPara4 = '<w:p w14:paraId="6C088A85" w14:textId="77777777" w:rsidR="007130E0" w:rsidRPr="005E7FD5" w:rsidRDefault="007130E0" w:rsidP="007130E0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Delete number of language support column in App </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>S</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">tore since we don’t have this </w:t></w:r></w:p>'

$para4.Range.InsertXML((New-PkgXml $newPara4))

Write-Output "Done applying edits."
Write-Output $d.Content.Text
